$wb = $excel.ActiveWorkbook

# --- Add "WorkCreation" sheet (new sheet always lands at index 1; rename + move to end) ---
$null = $wb.Worksheets.Add()
$wb.Worksheets.Item(1).Name = "WorkCreation"
$wb.Worksheets.Item("WorkCreation").Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$wsWorkCreation = $wb.Worksheets.Item("WorkCreation")
$wsWorkCreation.Range("B1").Value = "creator"
$wsWorkCreation.Range("C1").Value = "date"
$wsWorkCreation.Range("D1").Value = "work"
$wsWorkCreation.Range("B2").Value = "ss-person:minor-gordon"
$wsWorkCreation.Range("C2").NumberFormat = "yyyy-mm-dd"
$wsWorkCreation.Range("C2").Value = (Get-Date -Year 2023 -Month 4 -Day 7 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$wsWorkCreation.Range("D2").Value = "ss-work:test"

# --- Add "Work" sheet (new sheet always lands at index 1; rename + move to end) ---
$null = $wb.Worksheets.Add()
$wb.Worksheets.Item(1).Name = "Work"
$wb.Worksheets.Item("Work").Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$wsWork = $wb.Worksheets.Item("Work")
$wsWork.Range("A1").Value = "@id"
$wsWork.Range("B1").Value = "title"
$wsWork.Range("A2").Value = "ss-work:test"
$wsWork.Range("B2").Value = "Test work"
